{"js": "// Insert six new bullet-list paragraphs right after the paragraph that\n// contains \"Concepts Implications support / confidence thresholds...\"\n// and before the \"State: Type according property values...\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"Concepts Implications support / confidence thresholds.\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not locate anchor paragraph for insertion.\");\n}\n\nconst newParagraphTexts = [\n  \"Relationship Statements: Rules (LHS, RHS in a Context by a Concept Property) stated via Upper Aligned types, instances, attributes, values Statements, Kinds, CSPOs and Resources (Concepts / Implications).\",\n  \"Relationships / Rules Model. Model Aggregation / Inference.\",\n  \"Model (Upper Aligned CSPOs):\",\n  \"Rule / Concept (schema) example: (Context, Kind/Statement/Resource, Concept, Kind/Statement/Resource);\",\n  \"Implication (instances): Statements, Kinds, Resources. Parse instances as new Rules / Concepts.\",\n  \"Rule Aggregation: Each Statement is itself a Rule stating a single fact building a Concept by means of its Kinds relationships / CSPO inter Statements occurrences.\"\n];\n\n// Insert in order, each time directly after the previous insertion point so\n// the final order matches the source list top-to-bottom.\nlet insertAfter = anchor;\nfor (const text of newParagraphTexts) {\n  insertAfter = insertAfter.insertParagraph(text, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Insert six new bullet-list paragraphs right after the paragraph that\n# contains \"Concepts Implications support / confidence thresholds...\"\n# and before the \"State: Type according property values...\" paragraph.\n\n$d = $word.ActiveDocument\n\n$anchorIndex = 0\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -like \"*Concepts Implications support / confidence thresholds.*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq 0) {\n    throw \"Could not locate anchor paragraph for insertion.\"\n}\n\n$newParagraphTexts = @(\n    \"Relationship Statements: Rules (LHS, RHS in a Context by a Concept Property) stated via Upper Aligned types, instances, attributes, values Statements, Kinds, CSPOs and Resources (Concepts / Implications).\",\n    \"Relationships / Rules Model. Model Aggregation / Inference.\",\n    \"Model (Upper Aligned CSPOs):\",\n    \"Rule / Concept (schema) example: (Context, Kind/Statement/Resource, Concept, Kind/Statement/Resource);\",\n    \"Implication (instances): Statements, Kinds, Resources. Parse instances as new Rules / Concepts.\",\n    \"Rule Aggregation: Each Statement is itself a Rule stating a single fact building a Concept by means of its Kinds relationships / CSPO inter Statements occurrences.\"\n)\n\n$insertAfterIndex = $anchorIndex\nforeach ($text in $newParagraphTexts) {\n    $rng = $d.Paragraphs($insertAfterIndex).Range\n    $rng.InsertParagraphAfter()\n    $insertAfterIndex = $insertAfterIndex + 1\n    $d.Paragraphs($insertAfterIndex).Range.Text = $text\n}\n"}
